$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append the new data row (row 48) mirroring the existing layout:
# A=Name (ru), B=name (var id), C=lag, D=aggregation, E=stationarity,
# F=is_noise, G=is_sna, H=freq
# (write B before A so the new shared strings land in the same order as
# the saved workbook: "export_and_stocks" then the Russian label)
$ws.Cells.Item(48, 2).Value = "export_and_stocks"
$ws.Cells.Item(48, 1).Value = "Экспорт и инвестиции в запасы в постоянных ценах"
$ws.Cells.Item(48, 3).Value = 27
$ws.Cells.Item(48, 4).Value = "last"
$ws.Cells.Item(48, 5).Value = 2
$ws.Cells.Item(48, 6).Value = 0
$ws.Cells.Item(48, 7).Value = 0
$ws.Cells.Item(48, 8).Value = "q"

# Scroll the view down and move the selection to the row below the new data,
# matching the author's saved view state.
$ws.Application.ActiveWindow.ScrollRow = 35
$ws.Range("A49").Select()

# Restore the workbook window size recorded in the saved file.
$excel.ActiveWindow.Width = 23040
$excel.ActiveWindow.Height = 9192

# Update the absolute path recorded for the workbook (folder rename
# macro_forecast_russia -> macro_nowcast_russia).
$wb.Path = "C:\Users\migareev\Documents\macro_nowcast_russia\data"
